$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.926.34"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").Value = "2.254.43"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'301.64"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").Value = "'92.08"
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("D7").Value = "'0.531"
$ws.Range("E7").Value = "  +3.43%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.482"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("E10").Value = "  +9.55%  "
$ws.Range("D11").Value = "'32.25"
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").Value = "'6.66"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").Value = "2.603.58"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'14.08"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "2.261.55"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "'0.756"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").Value = "41.781.44"
$ws.Range("E19").Value = "  +4.99%  "
$ws.Range("D20").Value = "'12.11"
$ws.Range("E20").Value = "  +9.40%  "
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").Value = "'66.97"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'240.97"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "'2.55"
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("D28").Value = "'23.89"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("E29").Value = "  +12.74%  "
$ws.Range("D30").Value = "'9.65"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("D31").Value = "'158.76"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "'33.73"
$ws.Range("E32").Value = "  +5.73%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("D36").Value = "'3.04"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  +2.73%  "
$ws.Range("D38").Value = "'0.116"
$ws.Range("E38").Value = "  +3.45%  "
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("E40").Value = "  +7.43%  "
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("D43").Value = "2.048.20"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "'19.65"
$ws.Range("E44").Value = "  +8.20%  "
$ws.Range("D45").Value = "'0.0279"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("D48").Value = "'2.84"
$ws.Range("E48").Value = "  +4.66%  "
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").Value = "'1.13"
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("D51").Value = "'158.76"
$ws.Range("E51").Value = "  +5.12%  "

Write-Output "Updated cryptos list"
